$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fill in B3 and C3 ---
$ws.Range("B3").Value = "MikeS@78"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:MikeS@78")
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("C3").Value = "Invalid"

# --- Row 4: hyperlink created while cell still reads "MikeS@78"; the
#     link's cached display text stays "MikeS@78" even after the cell
#     text is retyped to "MikeS" further down ---
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:MikeS@78", "", "", "MikeS@78")
$ws.Range("B4").Style = "Hyperlink"

# --- Row 5 ---
$ws.Range("B5").Value = "MikeS@"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:MikeS@")
$ws.Range("B5").Style = "Hyperlink"

$ws.Range("A5").Value = "michael@fakemail.com"
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:michael@fakemail.com")
$ws.Range("A5").Style = "Hyperlink"

$ws.Range("C5").Value = "Invalid"

# --- finish row 4 ---
$ws.Range("A4").Value = "smith@fakemail.com"
$ws.Range("B4").Value = "MikeS"
$ws.Range("C4").Value = "Invalid"

# --- selection / view update ---
$ws.Range("D7").Select() | Out-Null
